# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (column E, rows 16-51) is reversed from descending
# (2003 .. 1703) to ascending (1703 .. 2003) order, reflecting the new/updated
# account-statement database. The "Valor Mora" column (F) keeps its 48000
# value for every period except the oldest one (period 1703), which carries
# the special 29509 value - previously attached to the last row (period
# 1703 was last, now it is first), so F16/F51 effectively swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ascending period list that now fills E16:E51 (was previously descending).
$periods = @(
    "1703","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# Valor Mora (column F): every row is 48000 except the row that now holds the
# oldest period (1703, row 16), which keeps the 29509 figure that used to sit
# on the last row (period 1703 was row 51 before the reorder).
for ($row = 16; $row -le 51; $row++) {
    $ws.Cells.Item($row, 6).Value = 48000
}
$ws.Cells.Item(16, 6).Value = 29509
